{"js": "// Update the 25 \"three-digit \u00d7 one-digit\" practice equations in the\n// single table of the document. Each of the 5 content rows (table rows\n// 0, 4, 9, 14, 19 \u2014 the blank rows in between are left untouched) has 5\n// cells, one equation per cell. We replace the old equation text with\n// the new one, cell by cell, addressed by (row, col) so the edit is\n// unambiguous regardless of any text collisions between old/new values.\n\n// [tableRowIndex, columnIndex, oldText, newText]\nconst edits = [\n  [0, 0, \"887\u00d73=\", \"112\u00d77=\"],\n  [0, 1, \"670\u00d76=\", \"151\u00d79=\"],\n  [0, 2, \"636\u00d78=\", \"153\u00d77=\"],\n  [0, 3, \"979\u00d78=\", \"474\u00d77=\"],\n  [0, 4, \"124\u00d79=\", \"784\u00d77=\"],\n\n  [4, 0, \"122\u00d73=\", \"202\u00d72=\"],\n  [4, 1, \"525\u00d73=\", \"105\u00d76=\"],\n  [4, 2, \"175\u00d73=\", \"525\u00d76=\"],\n  [4, 3, \"587\u00d76=\", \"302\u00d75=\"],\n  [4, 4, \"519\u00d76=\", \"979\u00d78=\"],\n\n  [9, 0, \"618\u00d75=\", \"562\u00d79=\"],\n  [9, 1, \"219\u00d76=\", \"311\u00d76=\"],\n  [9, 2, \"934\u00d73=\", \"329\u00d79=\"],\n  [9, 3, \"325\u00d72=\", \"720\u00d75=\"],\n  [9, 4, \"251\u00d76=\", \"622\u00d79=\"],\n\n  [14, 0, \"225\u00d73=\", \"755\u00d76=\"],\n  [14, 1, \"765\u00d76=\", \"613\u00d74=\"],\n  [14, 2, \"319\u00d77=\", \"580\u00d78=\"],\n  [14, 3, \"148\u00d74=\", \"163\u00d76=\"],\n  [14, 4, \"258\u00d78=\", \"166\u00d73=\"],\n\n  [19, 0, \"706\u00d76=\", \"627\u00d78=\"],\n  [19, 1, \"747\u00d74=\", \"314\u00d76=\"],\n  [19, 2, \"443\u00d75=\", \"594\u00d73=\"],\n  [19, 3, \"921\u00d79=\", \"455\u00d77=\"],\n  [19, 4, \"962\u00d72=\", \"528\u00d77=\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nfor (const [row, col, oldText, newText] of edits) {\n  const cell = table.getCell(row, col);\n  cell.load(\"value\");\n  await context.sync();\n\n  if (cell.value === oldText) {\n    cell.value = newText;\n  } else {\n    // Fall back to a scoped search-and-replace inside this specific cell\n    // in case the cell text doesn't match verbatim (defensive only).\n    const hits = cell.body.search(oldText, { matchCase: true });\n    hits.load(\"items\");\n    await context.sync();\n    for (let i = 0; i < hits.items.length; i++) {\n      hits.items[i].insertText(newText, \"Replace\");\n    }\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the 25 \"three-digit x one-digit\" practice equations in the\n# single table of the document. The table has 20 rows x 5 columns; only\n# rows 1, 5, 10, 15, 20 (1-based) hold equation text - the rows in\n# between are blank answer rows and are left untouched. Each cell in a\n# content row holds exactly one equation, so we address every edit by\n# its (row, column) position, which is unambiguous even though a couple\n# of the new values happen to equal other old values elsewhere in the\n# table.\n\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\n# Each entry: row, column, old text (sans end-of-cell marker), new text\n$edits = @(\n    @(1, 1, \"887\u00d73=\", \"112\u00d77=\"),\n    @(1, 2, \"670\u00d76=\", \"151\u00d79=\"),\n    @(1, 3, \"636\u00d78=\", \"153\u00d77=\"),\n    @(1, 4, \"979\u00d78=\", \"474\u00d77=\"),\n    @(1, 5, \"124\u00d79=\", \"784\u00d77=\"),\n\n    @(5, 1, \"122\u00d73=\", \"202\u00d72=\"),\n    @(5, 2, \"525\u00d73=\", \"105\u00d76=\"),\n    @(5, 3, \"175\u00d73=\", \"525\u00d76=\"),\n    @(5, 4, \"587\u00d76=\", \"302\u00d75=\"),\n    @(5, 5, \"519\u00d76=\", \"979\u00d78=\"),\n\n    @(10, 1, \"618\u00d75=\", \"562\u00d79=\"),\n    @(10, 2, \"219\u00d76=\", \"311\u00d76=\"),\n    @(10, 3, \"934\u00d73=\", \"329\u00d79=\"),\n    @(10, 4, \"325\u00d72=\", \"720\u00d75=\"),\n    @(10, 5, \"251\u00d76=\", \"622\u00d79=\"),\n\n    @(15, 1, \"225\u00d73=\", \"755\u00d76=\"),\n    @(15, 2, \"765\u00d76=\", \"613\u00d74=\"),\n    @(15, 3, \"319\u00d77=\", \"580\u00d78=\"),\n    @(15, 4, \"148\u00d74=\", \"163\u00d76=\"),\n    @(15, 5, \"258\u00d78=\", \"166\u00d73=\"),\n\n    @(20, 1, \"706\u00d76=\", \"627\u00d78=\"),\n    @(20, 2, \"747\u00d74=\", \"314\u00d76=\"),\n    @(20, 3, \"443\u00d75=\", \"594\u00d73=\"),\n    @(20, 4, \"921\u00d79=\", \"455\u00d77=\"),\n    @(20, 5, \"962\u00d72=\", \"528\u00d77=\")\n)\n\nforeach ($edit in $edits) {\n    $row = $edit[0]\n    $col = $edit[1]\n    $oldText = $edit[2]\n    $newText = $edit[3]\n\n    $cell = $table.Cell($row, $col)\n    $cellRange = $cell.Range\n    $current = $cellRange.Text.TrimEnd([char]13, [char]7)\n\n    if ($current -eq $oldText) {\n        $cellRange.Text = $newText\n    } else {\n        # Defensive fallback: scoped find/replace within this cell only.\n        $find = $cellRange.Find\n        $find.ClearFormatting()\n        $find.Text = $oldText\n        $find.Replacement.ClearFormatting()\n        $find.Replacement.Text = $newText\n        $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n    }\n}\n"}
